# Nudge the DER diagram shapes on slide 13 ("DER [R01]") up and to the
# left by a small, uniform amount (dx=-124178 EMU, dy=-36558 EMU), as done
# by selecting the diagram shapes (except the ITENSSAIDA box and one
# cardinality label) and nudging them in the PowerPoint UI.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

# Retângulo 5 (id=6)
$sh = $s.Shapes.Item(3)
$sh.Left = 261.3207092285156
$sh.Top = 317.9034118652344
# Retângulo 6 (id=7)
$sh = $s.Shapes.Item(4)
$sh.Left = 85.29252624511719
$sh.Top = 136.8896942138672
# Losango 8 (id=9)
$sh = $s.Shapes.Item(5)
$sh.Left = 120.78567504882812
$sh.Top = 327.38128662109375
# Conector de Seta Reta 17 (id=18)
$sh = $s.Shapes.Item(6)
$sh.Left = 182.64480590820312
$sh.Top = 352.8893127441406
# Conector reto 18 (id=19)
$sh = $s.Shapes.Item(7)
$sh.Left = 151.71511840820312
$sh.Top = 208.00244140625
# CaixaDeTexto 16 (id=17)
$sh = $s.Shapes.Item(8)
$sh.Left = 90.02488708496094
$sh.Top = 211.52276611328125
# CaixaDeTexto 19 (id=20)
$sh = $s.Shapes.Item(9)
$sh.Left = 211.250244140625
$sh.Top = 363.8567199707031
# Retângulo 10 (id=11)
$sh = $s.Shapes.Item(10)
$sh.Left = 567.05615234375
$sh.Top = 177.95559692382812
# Retângulo 11 (id=12)
$sh = $s.Shapes.Item(11)
$sh.Left = 581.2979736328125
$sh.Top = 449.79205322265625
# Retângulo 12 (id=13)
$sh = $s.Shapes.Item(12)
$sh.Left = 266.5492248535156
$sh.Top = 176.536865234375
# Retângulo 14 (id=15)
$sh = $s.Shapes.Item(13)
$sh.Left = 715.1326904296875
$sh.Top = 12.950394630432129
# Retângulo 20 (id=21)
$sh = $s.Shapes.Item(15)
$sh.Left = 744.6483764648438
$sh.Top = 308.425537109375
# Losango 21 (id=22)
$sh = $s.Shapes.Item(16)
$sh.Left = 325.8418273925781
$sh.Top = 401.1811828613281
# Losango 22 (id=23)
$sh = $s.Shapes.Item(17)
$sh.Left = 482.730712890625
$sh.Top = 468.747802734375
# Losango 23 (id=24)
$sh = $s.Shapes.Item(18)
$sh.Left = 801.3823852539062
$sh.Top = 439.071044921875
# Losango 24 (id=25)
$sh = $s.Shapes.Item(19)
$sh.Left = 790.2675170898438
$sh.Top = 182.49441528320312
# Losango 25 (id=26)
$sh = $s.Shapes.Item(20)
$sh.Left = 467.45111083984375
$sh.Top = 187.4334716796875
# Conector reto 26 (id=27)
$sh = $s.Shapes.Item(21)
$sh.Left = 450.5923767089844
$sh.Top = 213.64056396484375
# Conector reto 27 (id=28)
$sh = $s.Shapes.Item(22)
$sh.Left = 440.92041015625
$sh.Top = 494.2558288574219
# Conector reto 28 (id=29)
$sh = $s.Shapes.Item(23)
$sh.Left = 700.2662353515625
$sh.Top = 474.6740417480469
# Conector reto 29 (id=30)
$sh = $s.Shapes.Item(24)
$sh.Left = 830.6936645507812
$sh.Top = 340.7759094238281
# Conector reto 30 (id=31)
$sh = $s.Shapes.Item(25)
$sh.Left = 821.1971435546875
$sh.Top = 85.91803741455078
# Conector reto 31 (id=32)
$sh = $s.Shapes.Item(26)
$sh.Left = 695.6008911132812
$sh.Top = 209.1722869873047
# Losango 33 (id=34)
$sh = $s.Shapes.Item(27)
$sh.Left = 325.8418273925781
$sh.Top = 256.6979675292969
# Conector reto 34 (id=35)
$sh = $s.Shapes.Item(28)
$sh.Left = 355.3849792480469
$sh.Top = 211.52276611328125
# CaixaDeTexto 35 (id=36)
$sh = $s.Shapes.Item(29)
$sh.Left = 384.05181884765625
$sh.Top = 429.9325256347656
# CaixaDeTexto 36 (id=37)
$sh = $s.Shapes.Item(30)
$sh.Left = 381.02813720703125
$sh.Top = 387.0671081542969
# CaixaDeTexto 37 (id=38)
$sh = $s.Shapes.Item(31)
$sh.Left = 382.73095703125
$sh.Top = 294.556640625
# CaixaDeTexto 38 (id=39)
$sh = $s.Shapes.Item(32)
$sh.Left = 373.551513671875
$sh.Top = 249.8814239501953
# CaixaDeTexto 39 (id=40)
$sh = $s.Shapes.Item(33)
$sh.Left = 550.8346557617188
$sh.Top = 464.1784362792969
# CaixaDeTexto 40 (id=41)
$sh = $s.Shapes.Item(34)
$sh.Left = 450.5923767089844
$sh.Top = 467.9652099609375
# CaixaDeTexto 41 (id=42)
$sh = $s.Shapes.Item(35)
$sh.Left = 450.5923767089844
$sh.Top = 186.73977661132812
# CaixaDeTexto 42 (id=43)
$sh = $s.Shapes.Item(36)
$sh.Left = 540.5986938476562
$sh.Top = 189.89283752441406
# CaixaDeTexto 43 (id=44)
$sh = $s.Shapes.Item(37)
$sh.Left = 795.3778076171875
$sh.Top = 80.83787536621094
# CaixaDeTexto 44 (id=45)
$sh = $s.Shapes.Item(38)
$sh.Left = 736.51708984375
$sh.Top = 177.4735565185547
# CaixaDeTexto 46 (id=47)
$sh = $s.Shapes.Item(40)
$sh.Left = 754.6372680664062
$sh.Top = 477.73638916015625
